$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 34..50 get their full contents (columns A..AY) permuted: the data that
# currently lives in row $src ends up in row $dst (a pure rearrangement of the
# 17 existing rows - no values are altered).
$mapping = @{
    34 = 37
    35 = 44
    36 = 50
    37 = 34
    38 = 46
    39 = 35
    40 = 47
    41 = 38
    42 = 49
    43 = 41
    44 = 39
    45 = 48
    46 = 42
    47 = 36
    48 = 43
    49 = 40
    50 = 45
}

$firstCol = 1   # A
$lastCol  = 51  # AY

# Snapshot every source row's values before any writes happen, so overlapping
# source/destination rows don't clobber data we still need to read.
$snapshots = @{}
foreach ($dstRow in $mapping.Keys) {
    $srcRow = $mapping[$dstRow]
    if (-not $snapshots.ContainsKey($srcRow)) {
        $srcRange = $ws.Range($ws.Cells.Item($srcRow, $firstCol), $ws.Cells.Item($srcRow, $lastCol))
        $snapshots[$srcRow] = $srcRange.Value2
    }
}

foreach ($dstRow in $mapping.Keys) {
    $srcRow = $mapping[$dstRow]
    $dstRange = $ws.Range($ws.Cells.Item($dstRow, $firstCol), $ws.Cells.Item($dstRow, $lastCol))
    $dstRange.Value2 = $snapshots[$srcRow]
}
